# Update cryptocurrency price/volume data on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure Price (D) and Volume (E) columns keep their original text formatting
# so values like "1.011" are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @{
    2 = @("26.448.88", '  -3.49%  ')
    3 = @("1.801.99", '  -3.23%  ')
    4 = @("1.011", '  +0.70%  ')
    5 = @("1.010", '  +0.62%  ')
    6 = @("308.36", '  -2.12%  ')
    7 = @("0.4541", '  -1.87%  ')
    8 = @("0.3652", '  -1.79%  ')
    9 = @("0.07116", '  -2.79%  ')
    10 = @("0.8697", '  -2.31%  ')
    11 = @($null, '  -0.46%  ')
    12 = @("19.14", '  -4.44%  ')
    13 = @("1.838.82", '  -0.63%  ')
    14 = @("5.266", '  -2.47%  ')
    15 = @("6.321", '  -3.90%  ')
    16 = @("86.17", '  -6.33%  ')
    17 = @($null, '  +0.76%  ')
    18 = @("0.000008548", '  -4.60%  ')
    19 = @("1.009", '  +0.55%  ')
    20 = @("26.506.58", '  -3.32%  ')
    21 = @("14.21", '  -3.96%  ')
    22 = @("4.951", '  -3.57%  ')
    23 = @("2.062.66", '  -0.45%  ')
    24 = @("10.33", '  -2.13%  ')
    25 = @("1.984", '  +3.09%  ')
    26 = @("151.01", '  -0.88%  ')
    27 = @("17.81", '  -3.53%  ')
    28 = @("1.990", '  -3.41%  ')
    29 = @("112.53", '  -3.33%  ')
    30 = @("4.854", '  -4.93%  ')
    31 = @("0.08693", '  -1.81%  ')
    32 = @("3.074", '  -2.09%  ')
    33 = @("0.7279", '  -5.43%  ')
    34 = @("4.428", '  -1.96%  ')
    35 = @("1.109", '  -5.68%  ')
    36 = @("1.011", '  +1.00%  ')
    37 = @("2.501", '  -8.63%  ')
    38 = @($null, '  -0.65%  ')
    39 = @("0.01913", '  -2.50%  ')
    40 = @("0.05085", '  -3.12%  ')
    41 = @($null, '  -3.40%  ')
    42 = @("6.868", '  -2.93%  ')
    43 = @("0.4895", '  -4.97%  ')
    44 = @("0.1562", '  -4.88%  ')
    45 = @("8.118", '  -3.64%  ')
    46 = @("1.010", '  +0.66%  ')
    47 = @("0.4584", '  -4.87%  ')
    48 = @("101.71", '  -1.39%  ')
    49 = @("9.895", '  -4.15%  ')
    50 = @("1.578", '  -4.66%  ')
    51 = @("0.05989", '  -3.73%  ')
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
    $ws.Cells.Item($row, 5).Value = $eVal
}

$wb.Save()